$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: shift values to match old row3 (Target cluster -> FAPs) with recomputed TPM-based stats
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.026619
$ws.Range("H2").Value = 0.079857
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.016376
$ws.Range("N2").Value = 0.049128
$ws.Range("O2").Value = 0.4917717717717718
$ws.Range("P2").Value = 0.4917717717717718
$ws.Range("Q2").Value = 0.000435912744
$ws.Range("R2").Value = 0.003923214695999999
$ws.Range("S2").Value = 0.4917717717717718
$ws.Range("T2").Value = 0.4917717717717718

# Row 3: shift values to match old row4 (Target cluster -> MuSCs) with recomputed TPM-based stats
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.026619
$ws.Range("H3").Value = 0.079857
$ws.Range("M3").Value = 0.016924
$ws.Range("N3").Value = 0.050772
$ws.Range("O3").Value = 0.5082282282282282
$ws.Range("P3").Value = 0.5082282282282282
$ws.Range("Q3").Value = 0.000450499956
$ws.Range("R3").Value = 0.004054499604
$ws.Range("S3").Value = 0.5082282282282282
$ws.Range("T3").Value = 0.5082282282282282

# Row 4 no longer exists in the updated dataset - remove it entirely
$ws.Rows(4).Delete()
